# Update results in column AF (doctor_MA / std Dev) for rows 4-13
# as reflected in the latest stats run ("updated results and code").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    4  = 0.583
    5  = 0.857
    6  = 0.694
    7  = 0.783
    8  = 0.804
    9  = 0.714
    10 = 0.857
    11 = 0.857
    12 = 1.167
    13 = 2
}

foreach ($row in $updates.Keys) {
    $ws.Range("AF$row").Value = $updates[$row]
}
